$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 876.375
$ws.Range("I28").Value = 839.0833
$ws.Range("J28").Value = 988.25
$ws.Range("K28").Value = 839.0833
$ws.Range("L28").Value = 988.25
$ws.Range("M28").Value = -354.0833
$ws.Range("N28").Value = -1958.25
$ws.Range("H39").Value = 623.1111
$ws.Range("I39").Value = 658.2857
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 1974.8571
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = -1678.8571
$ws.Range("N39").Value = -2092
$ws.Range("H40").Value = 1515.3846
$ws.Range("I40").Value = 1580.1
$ws.Range("J40").Value = 1299.6666
$ws.Range("K40").Value = 1580.1
$ws.Range("L40").Value = 1299.6666
$ws.Range("M40").Value = -1405.1
$ws.Range("N40").Value = -1649.6666
$ws.Range("H98").Value = 287559.28
$ws.Range("I98").Value = 287559.28
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 287559.28
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -286061.28
$ws.Range("H122").Value = 287559.28
$ws.Range("I122").Value = 287559.28
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 862677.8400000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -860227.8400000001
$ws.Range("H132").Value = 1378454.1
$ws.Range("I132").Value = 1653550
$ws.Range("J132").Value = 2974.7273
$ws.Range("K132").Value = 4960650
$ws.Range("L132").Value = 8924.1819
$ws.Range("M132").Value = -4958120
$ws.Range("N132").Value = -13984.1819
$ws.Range("H137").Value = 1279.421
$ws.Range("I137").Value = 1017.1667
$ws.Range("J137").Value = 2262.875
$ws.Range("K137").Value = 3051.5001
$ws.Range("L137").Value = 6788.625
$ws.Range("M137").Value = -501.5001000000002
$ws.Range("N137").Value = -11888.625
$ws.Range("H138").Value = 3164.2317
$ws.Range("I138").Value = 636.569
$ws.Range("J138").Value = 7126.5137
$ws.Range("K138").Value = 1909.707
$ws.Range("L138").Value = 21379.5411
$ws.Range("M138").Value = 3230.293
$ws.Range("N138").Value = -31659.5411

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8421.215
$ws.Range("I32").Value = 4170.828
$ws.Range("J32").Value = 26556.2
$ws.Range("K32").Value = 4170.828
$ws.Range("L32").Value = 26556.2
$ws.Range("M32").Value = -3883.828
$ws.Range("N32").Value = -27130.2
$ws.Range("H45").Value = 71430590
$ws.Range("I45").Value = 83335464
$ws.Range("J45").Value = 1362
$ws.Range("K45").Value = 83335464
$ws.Range("L45").Value = 1362
$ws.Range("M45").Value = -83335087
$ws.Range("N45").Value = -2116
$ws.Range("H139").Value = 34043.832
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 34043.832
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 34043.832
$ws.Range("N139").Value = -44323.832

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1528.0476
$ws.Range("I134").Value = 1042.96
$ws.Range("J134").Value = 2241.4119
$ws.Range("K134").Value = 3128.88
$ws.Range("L134").Value = 6724.2357
$ws.Range("M134").Value = -593.8800000000001
$ws.Range("N134").Value = -11794.2357

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1638.6912
$ws.Range("I31").Value = 1066.3898
$ws.Range("J31").Value = 5390.4443
$ws.Range("K31").Value = 1066.3898
$ws.Range("L31").Value = 5390.4443
$ws.Range("M31").Value = -771.3897999999999
$ws.Range("N31").Value = -5980.4443
$ws.Range("H34").Value = 1638.6912
$ws.Range("I34").Value = 1066.3898
$ws.Range("J34").Value = 5390.4443
$ws.Range("K34").Value = 1066.3898
$ws.Range("L34").Value = 5390.4443
$ws.Range("M34").Value = -864.3897999999999
$ws.Range("N34").Value = -5794.4443
$ws.Range("H58").Value = 917.1905
$ws.Range("I58").Value = 724.48486
$ws.Range("J58").Value = 1623.7778
$ws.Range("K58").Value = 724.48486
$ws.Range("L58").Value = 1623.7778
$ws.Range("M58").Value = -521.48486
$ws.Range("N58").Value = -2029.7778
$ws.Range("H132").Value = 1222.3959
$ws.Range("I132").Value = 972.675
$ws.Range("J132").Value = 2471
$ws.Range("K132").Value = 2918.025
$ws.Range("L132").Value = 7413
$ws.Range("M132").Value = -388.0249999999996
$ws.Range("N132").Value = -12473
$ws.Range("H136").Value = 917.1905
$ws.Range("I136").Value = 724.48486
$ws.Range("J136").Value = 1623.7778
$ws.Range("K136").Value = 2173.45458
$ws.Range("L136").Value = 4871.3334
$ws.Range("M136").Value = 376.5454199999999
$ws.Range("N136").Value = -9971.3334

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H131").Value = 724.63635
$ws.Range("I131").Value = 453.0909
$ws.Range("J131").Value = 996.1818
$ws.Range("K131").Value = 1359.2727
$ws.Range("L131").Value = 2988.5454
$ws.Range("M131").Value = 3680.7273
$ws.Range("N131").Value = -13068.5454

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 20000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 20000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 20000
$ws.Range("N38").Value = -20926
$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 1909.3334
$ws.Range("J102").Value = 4090.6667
$ws.Range("K102").Value = 1909.3334
$ws.Range("L102").Value = 4090.6667
$ws.Range("M102").Value = -287.3334
$ws.Range("N102").Value = -7334.6667
$ws.Range("H126").Value = 2003.6364
$ws.Range("I126").Value = 1860
$ws.Range("J126").Value = 2123.3333
$ws.Range("K126").Value = 5580
$ws.Range("L126").Value = 6369.999899999999
$ws.Range("M126").Value = -3110
$ws.Range("N126").Value = -11309.9999

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2471.7144
$ws.Range("I40").Value = 2100.4
$ws.Range("J40").Value = 3400
$ws.Range("K40").Value = 2100.4
$ws.Range("L40").Value = 3400
$ws.Range("M40").Value = -1964.4
$ws.Range("N40").Value = -3672
$ws.Range("H122").Value = 2985.7632
$ws.Range("I122").Value = 2907.4707
$ws.Range("J122").Value = 3651.25
$ws.Range("K122").Value = 8722.4121
$ws.Range("L122").Value = 10953.75
$ws.Range("M122").Value = -6272.4121
$ws.Range("N122").Value = -15853.75
$ws.Range("H132").Value = 1094.6111
$ws.Range("I132").Value = 875.4792
$ws.Range("J132").Value = 2847.6667
$ws.Range("K132").Value = 2626.4376
$ws.Range("L132").Value = 8543.000100000001
$ws.Range("M132").Value = -96.4376000000002
$ws.Range("N132").Value = -13603.0001
$ws.Range("H136").Value = 1771.6666
$ws.Range("I136").Value = 938.5185
$ws.Range("J136").Value = 3271.3333
$ws.Range("K136").Value = 2815.5555
$ws.Range("L136").Value = 9813.999899999999
$ws.Range("M136").Value = -265.5554999999999
$ws.Range("N136").Value = -14913.9999

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5569.3
$ws.Range("I107").Value = 354.15384
$ws.Range("J107").Value = 15254.571
$ws.Range("K107").Value = 1062.46152
$ws.Range("L107").Value = 45763.713
$ws.Range("M107").Value = 857.5384799999999
$ws.Range("N107").Value = -49603.713
$ws.Range("H126").Value = 1445.6154
$ws.Range("I126").Value = 988.8
$ws.Range("J126").Value = 2968.3333
$ws.Range("K126").Value = 2966.4
$ws.Range("L126").Value = 8904.999899999999
$ws.Range("M126").Value = -496.3999999999996
$ws.Range("N126").Value = -13844.9999
$ws.Range("H128").Value = 39128.75
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 39128.75
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 39128.75
$ws.Range("N128").Value = -49088.75
